$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("E:E").Insert()
$ws.Columns("E:E").ColumnWidth = 24.42578125
Write-Host "done"
